{"js": "// Highlight specific \"ToDo\" bullet items with their review colors.\n// Yellow: items flagged for follow-up; Green: items already covered elsewhere.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst YELLOW_TEXTS = [\n  \"Improve frustum culling on stereo camera\",\n  \"Implement components update groups, everything inside an update group can be run in parallel (es. All rigid body in pyscis, or all colliders)\",\n  \"Integrate pysicx colliders, more efficient than c# ones\",\n  \"Think to move some material properties at shader level (es. Write depth, etc)\"\n];\n\nconst GREEN_TEXTS = [\n  \"Implement oculus depth map\"\n];\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const text = para.text;\n  if (YELLOW_TEXTS.indexOf(text) !== -1) {\n    para.font.highlightColor = \"Yellow\";\n  } else if (GREEN_TEXTS.indexOf(text) !== -1) {\n    para.font.highlightColor = \"BrightGreen\";\n  }\n}\n\nawait context.sync();\n", "ps1": "# Highlight specific \"ToDo\" bullet items with their review colors.\n# Yellow: items flagged for follow-up; Green: items already covered elsewhere.\n$d = $word.ActiveDocument\n\n$yellowTexts = @(\n    \"Improve frustum culling on stereo camera\",\n    \"Implement components update groups, everything inside an update group can be run in parallel (es. All rigid body in pyscis, or all colliders)\",\n    \"Integrate pysicx colliders, more efficient than c# ones\",\n    \"Think to move some material properties at shader level (es. Write depth, etc)\"\n)\n\n$greenTexts = @(\n    \"Implement oculus depth map\"\n)\n\nforeach ($p in $d.Paragraphs) {\n    $text = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($yellowTexts -contains $text) {\n        $p.Range.Select()\n        $word.Selection.Font.HighlightColorIndex = \"wdYellow\"\n    } elseif ($greenTexts -contains $text) {\n        $p.Range.Select()\n        $word.Selection.Font.HighlightColorIndex = \"wdBrightGreen\"\n    }\n}\n"}
